# fix: add point to frontend plan
#
# Locate the "Фронтенд" plans slide (shape "Объект 2") by scanning for the
# bullet that ends in "...http-запросы" and append a new level-2 bullet:
#   "Сделать сервисную страницу для администратора вендора"

$p = $ppt.ActivePresentation

$targetShape = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*http-запросы*") {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# The paragraph that currently ends the text body is the
# "...http-запросы" bullet; start a new paragraph right after it at the
# same (level-2) indent as the other bullets.
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count)

$part1 = $lastPara.InsertAfter("`rСделать сервисную страницу для ")
$lastPara2 = $tr.Paragraphs($tr.Paragraphs().Count)
$part2 = $lastPara2.InsertAfter("администратора вендора")
